$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("A16").Value = "ATA_WEAPON_TENTACLE_GUN"
$ws.Range("B16").Value = "Tentacle Gun"
$ws.Range("C16").Value = "触手枪"

$ws.Range("G9").Select()
